$wb = $excel.ActiveWorkbook

# --- MenuF sheet: remove two blank spacer rows (shifts the menu table up by 2)
# and add the new "Mains EU/US Select" sub-menu entries to the POWER IC CALIB MENU row.
$ws = $wb.Worksheets.Item("MenuF")
$ws.Rows("19:20").Delete()

$ws.Range("G26").Value = "Mains EU/US Select  "
$ws.Range("H26").Value = "  220V / 50 Hz  ?   "
$ws.Range("I26").Value = "  110V / 60 Hz  ?   "
$ws.Range("J26").Value = "  Mains Updated !   "

# Match the highlighted/boxed style used by the other "last cell in row" entries.
$ws.Range("F26").Copy()
$ws.Range("J26").PasteSpecial(-4122)

# --- The MenuF sheet becomes the active tab/selection (was Key before).
$ws.Activate()
$ws.Range("G26").Select()
